$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(17).Insert()
$ws.Range("A24:F24").Copy()
$ws.Range("A17").Select()
$ws.Paste()
$excel.CutCopyMode = $false
